# Hortaliza, Feria Lagunitas de Puerto Montt - Cebolla
# Weekly update: two new price records inserted at rows 457-458, pushing
# the previously existing rows 457-539 down to 459-541.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at 457:458 (existing data shifts down to 459:541).
$ws.Range("A457:R458").EntireRow.Insert()

# --- New row 457 ---
$ws.Range("A457").Value = 4
$ws.Range("B457").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C457").Value = "Los Lagos"
$ws.Range("D457").Value = 44694
$ws.Range("E457").Value = 10
$ws.Range("F457").Value = 100112004
$ws.Range("G457").Value = "Cebolla"
$ws.Range("H457").Value = "Morada(o)"
$ws.Range("I457").Value = "1a (cosecha)"
$ws.Range("J457").Value = 250
$ws.Range("K457").Value = 14000
$ws.Range("L457").Value = 14000
$ws.Range("M457").Value = 14000
$ws.Range("N457").Value = "`$/malla 18 kilos"
$ws.Range("O457").Value = "Región de O'Higgins"
$ws.Range("P457").Value = 778
$ws.Range("Q457").Value = 18
$ws.Range("R457").Value = "Hortaliza"

# --- New row 458 ---
$ws.Range("A458").Value = 4
$ws.Range("B458").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C458").Value = "Los Lagos"
$ws.Range("D458").Value = 44694
$ws.Range("E458").Value = 10
$ws.Range("F458").Value = 100112004
$ws.Range("G458").Value = "Cebolla"
$ws.Range("H458").Value = "Sin especificar"
$ws.Range("I458").Value = "1a (cosecha)"
$ws.Range("J458").Value = 900
$ws.Range("K458").Value = 8500
$ws.Range("L458").Value = 9000
$ws.Range("M458").Value = 8750
$ws.Range("N458").Value = "`$/malla 18 kilos"
$ws.Range("O458").Value = "Región de O'Higgins"
$ws.Range("P458").Value = 486
$ws.Range("Q458").Value = 18
$ws.Range("R458").Value = "Hortaliza"
